$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# --- Row 1 header: rename "EditDesc" -> "EditCategory" and give it the
#     same highlighted header look as the other key columns (B1/C1) ---
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "EditCategory"

# --- Row 2 sample data refreshed for the new GL Batch 1 run ---
$ws.Range("B2").Value = "AutomationTestGNB299"
$ws.Range("C2").Value = "AutomationTestGNB399"
$ws.Range("E2").Value = "Edit_Automation_GNB399"

# --- Row 4 helper text: add matching placeholder under the new column ---
$ws.Range("C4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "Provide unique AlphaNumeric or Numeric value here"

$excel.CutCopyMode = $false

# --- Restore the saved selection shown in the workbook (column shifted
#     from D10 to E10) ---
$ws.Range("E10").Select()
